$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 83333680
$ws.Range("I28").Value = 83333680
$ws.Range("K28").Value = 83333680
$ws.Range("M28").Value = -83333195
$ws.Range("H33").Value = 249.64285
$ws.Range("I33").Value = 226.81818
$ws.Range("J33").Value = 333.33334
$ws.Range("K33").Value = 226.81818
$ws.Range("L33").Value = 333.33334
$ws.Range("M33").Value = 2.181819999999988
$ws.Range("N33").Value = -791.33334
$ws.Range("H86").Value = 375501
$ws.Range("H89").Value = 375501
$ws.Range("H96").Value = 200002080
$ws.Range("I96").Value = 200002080
$ws.Range("K96").Value = 600006240
$ws.Range("M96").Value = -600004867
$ws.Range("H98").Value = 3275.25
$ws.Range("I98").Value = 3275.25
$ws.Range("K98").Value = 3275.25
$ws.Range("M98").Value = -1777.25
$ws.Range("H100").Value = 1642.4445
$ws.Range("I100").Value = 1635.25
$ws.Range("K100").Value = 1635.25
$ws.Range("M100").Value = -1094.25
$ws.Range("H103").Value = 3567.1667
$ws.Range("I103").Value = 3975.75
$ws.Range("J103").Value = 2750
$ws.Range("K103").Value = 11927.25
$ws.Range("L103").Value = 8250
$ws.Range("M103").Value = -11341.25
$ws.Range("N103").Value = -9422
$ws.Range("H122").Value = 3275.25
$ws.Range("I122").Value = 3275.25
$ws.Range("K122").Value = 9825.75
$ws.Range("M122").Value = -7375.75
$ws.Range("H132").Value = 2910.6667
$ws.Range("I132").Value = 1231.8096
$ws.Range("K132").Value = 3695.4288
$ws.Range("M132").Value = -1165.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16667859
$ws.Range("I2").Value = 19608546
$ws.Range("J2").Value = 3965.3333
$ws.Range("K2").Value = 19608546
$ws.Range("L2").Value = 3965.3333
$ws.Range("M2").Value = -19608433
$ws.Range("N2").Value = -4191.3333
$ws.Range("H38").Value = 5006.3335
$ws.Range("I38").Value = 5006.3335
$ws.Range("K38").Value = 5006.3335
$ws.Range("M38").Value = -4539.3335
$ws.Range("H61").Value = 5405.5
$ws.Range("I61").Value = 5859.4116
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 5859.4116
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -5647.4116
$ws.Range("N61").Value = -3257.3333
$ws.Range("H110").Value = 4275909.5
$ws.Range("I110").Value = 7408950
$ws.Range("J110").Value = 3581.818
$ws.Range("K110").Value = 7408950
$ws.Range("L110").Value = 3581.818
$ws.Range("M110").Value = -7406905
$ws.Range("N110").Value = -7671.818
$ws.Range("H116").Value = 16667859
$ws.Range("I116").Value = 19608546
$ws.Range("J116").Value = 3965.3333
$ws.Range("K116").Value = 19608546
$ws.Range("L116").Value = 3965.3333
$ws.Range("M116").Value = -19606252
$ws.Range("N116").Value = -8553.3333
$ws.Range("H136").Value = 5405.5
$ws.Range("I136").Value = 5859.4116
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 17578.2348
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -15028.2348
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16667859
$ws.Range("I3").Value = 19608546
$ws.Range("J3").Value = 3965.3333
$ws.Range("K3").Value = 19608546
$ws.Range("L3").Value = 3965.3333
$ws.Range("M3").Value = -19608432
$ws.Range("N3").Value = -4193.3333
$ws.Range("H105").Value = 3972366
$ws.Range("I105").Value = 5955891
$ws.Range("K105").Value = 5955891
$ws.Range("M105").Value = -5954144
$ws.Range("H107").Value = 2866.6667
$ws.Range("I107").Value = 2866.6667
$ws.Range("K107").Value = 2866.6667
$ws.Range("M107").Value = -946.6667000000002
$ws.Range("H134").Value = 2180
$ws.Range("I134").Value = 2440.2222
$ws.Range("K134").Value = 7320.6666
$ws.Range("M134").Value = -4785.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2938.487
$ws.Range("I31").Value = 2223.9375
$ws.Range("K31").Value = 2223.9375
$ws.Range("M31").Value = -1928.9375
$ws.Range("H34").Value = 2938.487
$ws.Range("I34").Value = 2223.9375
$ws.Range("K34").Value = 2223.9375
$ws.Range("M34").Value = -2021.9375
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24127500
$ws.Range("I4").Value = 37799544
$ws.Range("J4").Value = 366.41177
$ws.Range("K4").Value = 113398632
$ws.Range("L4").Value = 1099.23531
$ws.Range("M4").Value = -113398520
$ws.Range("N4").Value = -1323.23531
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 250
$ws.Range("K19").Value = 250
$ws.Range("M19").Value = 38

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5006750
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 3329.8462
$ws.Range("I7").Value = 3329.8462
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3329.8462
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3217.8462
$ws.Range("N7").ClearContents()
$ws.Range("H55").Value = 624.4211
$ws.Range("I55").Value = 586.1667
$ws.Range("J55").Value = 690
$ws.Range("K55").Value = 586.1667
$ws.Range("L55").Value = 690
$ws.Range("M55").Value = -413.1667
$ws.Range("N55").Value = -1036
$ws.Range("H82").Value = 63852.938
$ws.Range("I82").Value = 2499.5
$ws.Range("J82").Value = 84304.086
$ws.Range("K82").Value = 2499.5
$ws.Range("L82").Value = 84304.086
$ws.Range("M82").Value = -2138.5
$ws.Range("N82").Value = -85026.086
$ws.Range("H85").Value = 63852.938
$ws.Range("I85").Value = 2499.5
$ws.Range("J85").Value = 84304.086
$ws.Range("K85").Value = 2499.5
$ws.Range("L85").Value = 84304.086
$ws.Range("M85").Value = -1251.5
$ws.Range("N85").Value = -86800.086
$ws.Range("H126").Value = 3329.8462
$ws.Range("I126").Value = 3329.8462
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9989.5386
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7519.5386
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 6515.909
$ws.Range("I136").Value = 5644
$ws.Range("J136").Value = 8041.75
$ws.Range("K136").Value = 16932
$ws.Range("L136").Value = 24125.25
$ws.Range("M136").Value = -14382
$ws.Range("N136").Value = -29225.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 500412.5
$ws.Range("J2").Value = 550
$ws.Range("L2").Value = 550
$ws.Range("N2").Value = -774
$ws.Range("H82").Value = 28431.666
$ws.Range("J82").Value = 28431.666
$ws.Range("L82").Value = 28431.666
$ws.Range("N82").Value = -29197.666
$ws.Range("H85").Value = 28431.666
$ws.Range("J85").Value = 28431.666
$ws.Range("L85").Value = 28431.666
$ws.Range("N85").Value = -31083.666
$ws.Range("H96").Value = 1924.5
$ws.Range("I96").Value = 1924.5
$ws.Range("K96").Value = 1924.5
$ws.Range("M96").Value = -551.5
$ws.Range("H100").Value = 1360.2174
$ws.Range("I100").Value = 1590.3529
$ws.Range("J100").Value = 708.1667
$ws.Range("K100").Value = 3180.7058
$ws.Range("L100").Value = 1416.3334
$ws.Range("M100").Value = -2639.7058
$ws.Range("N100").Value = -2498.3334

Write-Output "done"